$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.549.09"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").Value = "3.021.54"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'595.86"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").Value = "'150.59"
$ws.Range("E6").Value = "  +6.86%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.017.56"
$ws.Range("E8").Value = "  +1.59%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  +10.99%  "

$ws.Range("E11").Value = "  +5.43%  "

$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("E13").Value = "  +3.71%  "

$ws.Range("D14").Value = "'34.60"
$ws.Range("E14").Value = "  +2.11%  "

$ws.Range("E15").Value = "  +2.63%  "

$ws.Range("D16").Value = "3.520.37"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").Value = "62.502.29"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").Value = "3.021.34"
$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("D20").Value = "'450.67"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "'14.19"
$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("D22").Value = "'0.690"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("D23").Value = "'7.46"
$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D24").Value = "'82.38"
$ws.Range("E24").Value = "  +1.60%  "

$ws.Range("D25").Value = "'10.93"
$ws.Range("E25").Value = "  +11.52%  "

$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +5.09%  "

$ws.Range("D27").Value = "'12.07"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").Value = "'2.72"
$ws.Range("E29").Value = "  +3.39%  "

$ws.Range("D30").Value = "'7.34"
$ws.Range("E30").Value = "  +7.72%  "

$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("E32").Value = "  +4.62%  "

$ws.Range("D33").Value = "'27.56"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("E34").Value = "  +3.10%  "

$ws.Range("D35").Value = "0.0₃0855"
$ws.Range("E35").Value = "  +11.01%  "

$ws.Range("E36").Value = "  +1.79%  "

$ws.Range("D37").Value = "'5.85"
$ws.Range("E37").Value = "  +3.11%  "

$ws.Range("D38").Value = "'3.06"
$ws.Range("E38").Value = "  +10.37%  "

$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("D40").Value = "'50.11"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "'9.05"
$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("E43").Value = "  +9.34%  "

$ws.Range("D44").Value = "'40.31"
$ws.Range("E44").Value = "  +9.20%  "

$ws.Range("D45").Value = "'391.02"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("D47").Value = "2.736.83"
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D48").Value = "'132.83"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").Value = "'2.19"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("E51").Value = "  +0.07%  "
